$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-147 down to 41-148.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new price-report record.
$ws.Cells.Item(40, 1).Value  = 7
$ws.Cells.Item(40, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(40, 3).Value  = "Ñuble"
$ws.Cells.Item(40, 4).Value  = 44414
$ws.Cells.Item(40, 5).Value  = 16
$ws.Cells.Item(40, 6).Value  = 100114013
$ws.Cells.Item(40, 7).Value  = "Zanahoria"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 120
$ws.Cells.Item(40, 11).Value = 5000
$ws.Cells.Item(40, 12).Value = 5500
$ws.Cells.Item(40, 13).Value = 5250
$ws.Cells.Item(40, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(40, 16).Value = 262
$ws.Cells.Item(40, 17).Value = 20
$ws.Cells.Item(40, 18).Value = "Hortaliza"
